$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Plain($addr, $text) {
    $ws.Range($addr).Value2 = $text
}

function Set-TextLiteral($addr, $text) {
    # Prevent Excel's auto date/number detection (e.g. "2025-11-01") from
    # turning the literal string into a date serial value; keep it as text.
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value2 = $text
    $c.Style = "Normal"
}

# --- Row 6 additions (Roteiro / Video Pronto / Data) ---
Set-Plain "D6" "E se eu te dissesse que o caos no Rio de Janeiro já foi previsto, em detalhes, anos atrás? Preste atenção.`nHoje, a cidade maravilhosa vive um pesadelo. Um clima de guerra. Megaoperações policiais se tornaram rotina, e a tensão nas ruas é quase palpável. Mas o que pouca gente sabe, é que esse cenário apocalíptico foi alertado. E não por qualquer um.`nVídeos antigos, de lives no YouTube, ressurgiram, mostrando um profeta descrevendo exatamente o que vemos agora. Ele falava de uma 'megaoperação' que pararia o Rio, de um 'clima de guerra' que se instalaria na capital fluminense. De uma cidade em 'caos declarado'. Pense nos combates diários, nos desabafos desesperados de figuras públicas. Nas famílias vivendo sob fogo cruzado. Não é assustador como tudo se encaixa?`nSeria mera coincidência? Ou estamos testemunhando algo que transcende nossa compreensão, onde o sobrenatural se choca com a dura realidade social? Essa conexão é inegável, e aterrorizante.`nAgora me diz, você acredita nessa profecia? Deixe sua opinião e compartilhe esse vídeo. O mundo precisa saber disso."
Set-Plain "E6" "OK"
Set-TextLiteral "G6" "2025-11-01"

# --- Row 7 additions (Roteiro / Video Pronto / Data) ---
Set-Plain "D7" "O lendário 'gato fantasma' da Espanha? Ele é real e foi fotografado!`nPor séculos, histórias sussurravam sobre um felino espectral, uma sombra branca que vagava pelas montanhas da Península Ibérica. Ninguém nunca teve certeza… até agora. Em Jaén, na Espanha, um fotógrafo, numa noite que parecia como qualquer outra, registrou o impensável. Uma visão que parou o tempo e fez o mundo prender a respiração.`nDiante da câmera, ali estava ele: um lince-ibérico branco. Não albino, mas leucístico. Seus olhos, vibrantes e selvagens, contrastavam com a pelagem pura como a neve. Um lince! Uma espécie que mal sobreviveu ao fio da navalha da extinção, agora revela um exemplar tão único que reescreve a própria história. Você consegue imaginar a emoção desse momento? A beleza é chocante. Mas com essa raridade vem uma preocupação: essa pelagem, embora deslumbrante, compromete sua camuflagem natural. Um fantasma que, paradoxalmente, se tornou mais visível.`nEsse felino fantasma agora é real. E sua existência levanta uma pergunta assustadora: sua beleza será sua salvação ou sua perdição? Comente o que você pensa e compartilhe essa história inacreditável."
Set-Plain "E7" "OK"
Set-TextLiteral "G7" "2025-11-01"

# --- Row 8 (new): Mistérios Cósmicos e Conspirações Alienígenas ---
Set-Plain "A8" "Mistérios Cósmicos e Conspirações Alienígenas"
Set-Plain "B8" "O que a NASA realmente esconde sobre o Cometa 3I/ATLAS? Com mais de 200 milhões de visualizações em diversos vídeos, este 'cometa' é teorizado como uma nave alienígena por anomalias em sua trajetória, composição e comportamento ao se aproximar do Sol. A especulação sobre 'imagens secretas' do Telescópio James Webb e a insistência de cientistas renomados como Avi Loeb em sua natureza extraterrestre criam um suspense fascinante. Este tema alimenta a curiosidade, o medo do desconhecido e a busca por respostas sobre vida inteligente no universo, garantindo engajamento viral."
Set-Plain "C8" "alta"
Set-Plain "D8" "Você sabia que o cometa mais misterioso da história pode não ser um cometa, mas sim… uma nave alienígena?`nPouca gente sabe disso, mas o que vou te contar sobre o Cometa 3I/ATLAS é real... e pode mudar sua visão do universo. O ATLAS, com mais de 200 milhões de visualizações em vídeos, virou um enigma cósmico que a NASA insiste em chamar de rocha. Mas a verdade é que sua trajetória era estranha. Sua composição, mais ainda. Ao se aproximar do Sol, ele não se comportava como um cometa normal. Era como se estivesse... manobrando.`nCientistas renomados, como Avi Loeb de Harvard, começaram a levantar a hipótese mais ousada: o ATLAS seria tecnologia extraterrestre. Uma nave disfarçada, observando a Terra. E o que a NASA realmente esconde? Há rumores insistentes de imagens secretas do Telescópio James Webb, que mostram algo além da nossa compreensão. Algo que foi rapidamente classificado, mas que levanta a pergunta: o que eles não querem que a gente veja?`nEnquanto o mundo dormia, os dados chegavam. As anomalias se acumulavam. Ninguém entendeu de fato o que estava acontecendo nos bastidores. Será que não estamos sozinhos? E se essa for a prova que sempre procuramos, escondida à vista de todos, por medo do que a revelação causaria?`nIsso muda tudo, não é? Agora me diz: você acredita que o Cometa ATLAS é uma nave alienígena? Comente sua teoria e mande esse vídeo pra quem precisa saber disso."
Set-Plain "E8" "OK"
Set-TextLiteral "G8" "2025-11-02"

# --- Row 9 (new): Heidi Klum e a Reinvenção do Halloween ---
Set-Plain "A9" "Heidi Klum e a Reinvenção do Halloween"
Set-Plain "B9" "Ninguém domina o Halloween como Heidi Klum, cujos disfarces anuais são um fenômeno global que paralisam a internet. Cada ano, como em 2025 com sua Medusa hiper-realista e seu marido Tom Kaulitz transformado em pedra, seus trajes e a festa 'Heidiween' em Nova Iorque levam meses de planejamento e geram um misto de fascínio, choque e antecipação. A viralidade vem da transformação extrema, dos detalhes chocantes e da curiosidade em ver o que a 'rainha do Halloween' fará a seguir, provando que criatividade e espetáculo são receita para o sucesso."
Set-Plain "C9" "alta"
Set-Plain "D9" "Pouca gente entende a obsessão de Heidi Klum pelo Halloween. Mas o que eu vou te contar agora, vai mudar sua visão sobre a rainha do terror.`nEnquanto o mundo se prepara para uma noite de sustos e doces, Heidi já planeja há meses. A verdade é que ninguém, absolutamente NINGUÉM, domina o Halloween como ela. É um fenômeno global que paralisa a internet a cada ano, transformando Nova Iorque em seu palco sombrio.`nVocê sabia que para 2025, os sussurros já apontam para uma Medusa hiper-realista, com escamas cintilantes e cobras que parecem rastejar? E seu marido, Tom Kaulitz, transformado em pedra, ao seu lado. É mais que uma fantasia; é uma obra de arte viva, uma cena de horror que ganha vida.`nMeses de trabalho, equipes de especialistas, detalhes tão chocantes que geram fascínio e repulsa. A 'Heidiween' não é só uma festa; é um ritual, uma antecipação do que a mente de Heidi Klum pode criar para nos chocar, nos intrigar. Você consegue imaginar o que é preciso para alcançar esse nível de transformação extrema?`nÉ essa mistura de criatividade insana e espetáculo puro que a torna a verdadeira rainha do Halloween. E cada ano, a pergunta ecoa: o que vem depois?`nAgora me diz... você teria coragem de encarar uma Medusa assim de perto? Manda esse vídeo pra quem precisa entender o verdadeiro poder de Heidi Klum."
Set-Plain "E9" "OK"
Set-TextLiteral "G9" "2025-11-02"

# --- Row 10 (new): Confrontos e Rivalidades de Personalidades ---
Set-Plain "A10" "Confrontos e Rivalidades de Personalidades"
Set-Plain "B10" "O Brasil parou para assistir a uma rivalidade explosiva que transcendeu as batalhas de rima e foi parar no ringue! A luta entre os MCs Xamuel e Jhony MC no Fight Music Show 7, que aconteceu em 01/11, gerou uma comoção nacional com milhões de visualizações em vídeos de prévias, encaradas tensas e reações pós-luta. A rivalidade acirrada, as provocações e a expectativa do confronto físico criaram um drama intenso e uma curiosidade imensa, alimentando debates sobre quem venceria e gerando conteúdos de reação em massa. Vídeos mostrando a encarada, falas polêmicas dos MCs e análises dos resultados são perfeitos para engajar a audiência com a tensão da disputa."
Set-Plain "C10" "alta"
Set-Plain "D10" "Pouca gente sabe disso, mas uma rivalidade explosiva que parecia ficção parou o Brasil no dia 01 de novembro. O que começou em batalhas de rima, onde as palavras eram as armas, escalou para um nível que ninguém imaginava. Estamos falando de Xamuel e Jhony MC. Dois nomes que, por meses, trocaram farpas e promessas, alimentando uma tensão que era quase palpável nas redes sociais. Cada provocação, cada encarada tensa antes da luta... parecia um capítulo de um suspense. E então, no Fight Music Show 7, a contagem regressiva para o confronto físico chegou ao fim. As luzes do ringue focavam neles, enquanto milhões de olhos grudados nas telas esperavam o desfecho. O ar estava pesado, a expectativa era imensa. Quem venceria? Essa era a pergunta que ecoava em cada casa, em cada comentário. E o que aconteceu lá dentro do ringue? Bom, o resultado não só definiu um vencedor, mas reacendeu toda a chama dessa rivalidade, gerando uma onda de reações e análises que inundou a internet. Essa tensão parou um país inteiro. Mas e você, o que pensa sobre o papel da rivalidade nesses palcos? Será que é só espetáculo ou existe algo mais? Comente sua opinião e compartilhe essa história com seus amigos!"
Set-Plain "E10" "OK"
Set-TextLiteral "G10" "2025-11-02"

# --- Row 11 (new): Crises Urbanas e Segurança Pública ---
Set-Plain "A11" "Crises Urbanas e Segurança Pública"
Set-Plain "B11" "Uma onda de medo e caos tomou conta do Rio de Janeiro, com operações policiais intensas que deixaram a população em alerta máximo! A recente megaoperação policial na cidade resultou em confrontos e gerou milhões de visualizações em vídeos de notícias, depoimentos de moradores e reações de personalidades como o ator André D'Lucca. A tensão das ruas, os relatos emocionantes de quem vive a realidade local, a presença policial e as discussões sobre violência e segurança pública criam um forte apelo emocional e um senso de urgência. Conteúdos que mostram a realidade da cidade, análises do impacto social e mensagens de solidariedade podem gerar grande engajamento e debate."
Set-Plain "C11" "alta"

# --- Row 12 (new): Transformações Impactantes de Celebridades ---
Set-Plain "A12" "Transformações Impactantes de Celebridades"
Set-Plain "B12" "Chocante! Heidi Klum, a rainha do Halloween, superou todas as expectativas com sua fantasia de Medusa, conquistando o título de melhor do mundo mais uma vez! Sua transformação inacreditável em Medusa para a tradicional festa 'Heidiween' em Nova York gerou mais de 558 milhões de visualizações, com detalhes impressionantes da maquiagem e do figurino. A curiosidade em torno da revelação, o impacto visual da criatividade e o suspense da transformação cativam o público, que busca ver o 'antes e depois' e os bastidores da criação. Vídeos que exploram o processo de elaboração da fantasia, a reação das pessoas e comparativos com anos anteriores geram fascínio e são altamente compartilháveis."
Set-Plain "C12" "alta"

# Move/resize selection to match the editor's final cursor position
$null = $ws.Range("A10:XFD11").Select()
